# "Upgrade to MVC 4 / EF5" -- Master Backlog maintenance pass on the
# "Generic Backlog" sheet: swap out completed / stale backlog items for
# the next batch of work, re-balance the "next up" (green) highlight
# band, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generic Backlog")

# ---------------------------------------------------------------------
# Color / style constants (direct formatting, matching the existing
# palette already used on this sheet).
# ---------------------------------------------------------------------
$YELLOW = 65535      # RGB(255,255,0)
$GREEN  = 5296274    # RGB(146,208,80)
$BLUE   = 15773696   # RGB(0,176,240)
$xlRight   = -4152
$xlPasteFormats = -4122

function Set-ItemStyle($rng, $bold, $color) {
    $rng.Font.Bold = $bold
    if ($color -eq $null) {
        $rng.Interior.ColorIndex = -4142
    } else {
        $rng.Interior.Color = $color
    }
}

# ---------------------------------------------------------------------
# 1) "Essential" section -- reword the massage/chiro backlog item.
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Repair the Machine = Ronnie + Acupuncture + PT"

# ---------------------------------------------------------------------
# 2) "Professional" section (rows 8-13) shrinks from 6 items to 5.
#    Re-purpose rows 8-12 with the new task names/status, then delete
#    the now-redundant row 13 (its content -- Get Web Root Antivirus --
#    slides up into row 12).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Push Market (Phase 2)"
$ws.Range("B8").Value = "IN PROGRESS"

$ws.Range("A9").Value = "Marketing Escalation for Terminus 2.0"

$ws.Range("A10").Value = "Office + Machine 2.0"

# Invoice Michael stays, but is no longer mid-list -- it now gets the
# same "next up" green highlight as the new last row.
$ws.Range("A11").Value = "Invoice Michael"
Set-ItemStyle $ws.Range("A11") $false $GREEN
Set-ItemStyle $ws.Range("B11") $false $GREEN
$ws.Range("B11").HorizontalAlignment = $xlRight

$ws.Range("A12").Value = "Get Web Root Antivirus"
Set-ItemStyle $ws.Range("A12") $false $GREEN
Set-ItemStyle $ws.Range("B12") $false $GREEN
$ws.Range("B12").HorizontalAlignment = $xlRight

$ws.Rows(13).Delete()

# ---------------------------------------------------------------------
# 3) "Personal / Household" section (now starting at row 14) grows
#    from 2 items to 3: "Contact Fran" becomes "Haircut", and a brand
#    new "Replacement Phone" row is inserted right after it (the old
#    Professional-section "Replacement Phone" task moved here).
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Haircut"

$ws.Rows(16).Insert()
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial($xlPasteFormats)
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("A16").Value = "Replacement Phone"
$ws.Range("B16").Value = "IN PROGRESS"

# "Use YNAB" is no longer the lone "next up" row, so it loses its bold
# and shares the same plain-green treatment as the rest of the band.
Set-ItemStyle $ws.Range("A17") $false $GREEN
Set-ItemStyle $ws.Range("B17") $false $GREEN
$ws.Range("B17").HorizontalAlignment = $xlRight

# ---------------------------------------------------------------------
# 4) "Other Stuff" section -- reword one item, and drop the bold on the
#    last row's blue highlight to match the rest of the sheet.
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "New T-Shirts and Shoes"

Set-ItemStyle $ws.Range("A22") $false $BLUE
Set-ItemStyle $ws.Range("B22") $false $BLUE
$ws.Range("B22").HorizontalAlignment = $xlRight

# ---------------------------------------------------------------------
# 5) Move the active selection.
# ---------------------------------------------------------------------
$ws.Range("A5").Select()
